$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (they store formatted strings,
# e.g. "26.425.85", "  +0.57%  ") rather than being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.425.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.696.99'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.58'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5489'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.43%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2734'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06446'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.97'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07674'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.714.92'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.557'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5851'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008411'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.70'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.485.80'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.946'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.57'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.259'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.89'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1314'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.918'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.81'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06248'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.384'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.331'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.617'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.597'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.691'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.040'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6171'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.410'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.760'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01652'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.120.25'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.120'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8809'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.24'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.848.64'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.65'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000108'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.208'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05288'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.126'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4303'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.05%  '
